$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each entry: row number followed by the 20 column values (A..T) in order
$rowsData = @(
    ,@(2,"Neutro","Agrp","Mc4r","sCs",1,0.3333333333333333,0.2239496666666667,0.671849,0.04591675498971699,0.04591675498971699,3,1,0.3318626666666666,0.9955879999999999,1,1,0.07432053357911111,0.668884802212,0.04591675498971699,0.04591675498971699)
    ,@(3,"FAPs","Agrp","Mc4r","sCs",2,0.6666666666666666,0.7402036666666666,2.220611,0.1517651305791486,0.1517651305791486,3,1,0.3318626666666666,0.9955879999999999,1,1,0.2456459626964444,2.210813664268,0.1517651305791486,0.1517651305791486)
    ,@(4,"M1","Agrp","Mc4r","sCs",1,0.3333333333333333,0.1166683333333333,0.350005,0.02392069323639076,0.02392069323639076,3,1,0.3318626666666666,0.9955879999999999,1,1,0.03871786421555555,0.34846077794,0.02392069323639076,0.02392069323639076)
    ,@(5,"M2","Agrp","Mc4r","sCs",1,0.3333333333333333,0.5808623333333334,1.742587,0.119095124540285,0.119095124540285,3,1,0.3318626666666666,0.9955879999999999,1,1,0.1927665229062222,1.734898706156,0.119095124540285,0.119095124540285)
    ,@(6,"Neutro","Agrp","Mc4r","sCs",3,1,3.024446,9.073338,0.6201069554094576,0.6201069554094577,3,1,0.3318626666666666,0.9955879999999999,1,1,1.003700714749333,9.033306432743998,0.6201069554094576,0.6201069554094577)
    ,@(7,"sCs","Agrp","Mc4r","sCs",1,0.3333333333333333,0.1911673333333333,0.573502,0.03919534124500099,0.03919534124500099,3,1,0.3318626666666666,0.9955879999999999,1,1,0.06344130101955554,0.5709717091759999,0.03919534124500099,0.03919534124500099)
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $entry[$i + 1]
    }
}